$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.032.47'
$ws.Range('E2').Value = '  -2.17%  '
$ws.Range('D3').Value = '3.583.30'
$ws.Range('E3').Value = '  -3.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.51'
$ws.Range('E5').Value = '  -6.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '191.95'
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('D7').Value = '3.579.63'
$ws.Range('E7').Value = '  -3.38%  '
$ws.Range('E8').Value = '  -3.02%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.677'
$ws.Range('E10').Value = '  -6.87%  '
$ws.Range('E11').Value = '  -6.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.51'
$ws.Range('E12').Value = '  -9.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000268'
$ws.Range('E13').Value = '  -6.68%  '
$ws.Range('E14').Value = '  -5.95%  '
$ws.Range('D15').Value = '4.150.09'
$ws.Range('E15').Value = '  -3.42%  '
$ws.Range('D16').Value = '3.582.39'
$ws.Range('E16').Value = '  -3.40%  '
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.37'
$ws.Range('E18').Value = '  -5.77%  '
$ws.Range('D19').Value = '66.919.48'
$ws.Range('E19').Value = '  -2.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.14'
$ws.Range('E20').Value = '  -5.61%  '
$ws.Range('E21').Value = '  -7.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '397.95'
$ws.Range('E22').Value = '  -2.83%  '
$ws.Range('E23').Value = '  -9.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.71'
$ws.Range('E24').Value = '  -4.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.23'
$ws.Range('E25').Value = '  -3.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.94'
$ws.Range('E26').Value = '  -4.27%  '
$ws.Range('E27').Value = '  -5.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.09'
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.93'
$ws.Range('E30').Value = '  -7.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.65'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.15'
$ws.Range('E32').Value = '  -5.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '630.20'
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.17'
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('E35').Value = '  -6.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '63.79'
$ws.Range('E36').Value = '  -5.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '41.97'
$ws.Range('E37').Value = '  -13.38%  '
$ws.Range('E38').Value = '  -3.97%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').Value = '0.0₃0760'
$ws.Range('E40').Value = '  -6.91%  '
$ws.Range('D41').Value = '3.152.89'
$ws.Range('E41').Value = '  +7.34%  '
$ws.Range('E42').Value = '  -4.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.71'
$ws.Range('E44').Value = '  +3.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.96'
$ws.Range('E45').Value = '  -2.98%  '
$ws.Range('E46').Value = '  -7.04%  '
$ws.Range('E47').Value = '  +2.61%  '
$ws.Range('E48').Value = '  -7.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.29'
$ws.Range('E49').Value = '  -4.32%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.51'
$ws.Range('E50').Value = '  -9.76%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.54'
$ws.Range('E51').Value = '  -2.62%  '
